$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.449548482894897
$ws.Range("B1").Value = 3.819766998291016
$ws.Range("C1").Value = 2.821645498275757
$ws.Range("D1").Value = 2.418932437896729
$ws.Range("E1").Value = 1.905972719192505
